$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing data rows 2-8 with new Transaccion / Cuotas values
$ws.Range("A2").Value = 63582208
$ws.Range("B2").Value = 1

$ws.Range("A3").Value = 948403273
$ws.Range("B3").Value = 1

$ws.Range("A4").Value = 530936366
$ws.Range("B4").Value = 3

$ws.Range("A5").Value = 852350050
$ws.Range("B5").Value = 6

$ws.Range("A6").Value = 690758213
$ws.Range("B6").Value = 1

$ws.Range("A7").Value = 517999539
$ws.Range("B7").Value = 3

$ws.Range("A8").Value = 504974250
$ws.Range("B8").Value = 6

# Remove the now-unused rows 9-14 entirely so the sheet dimension shrinks to A1:B8
$xlShiftUp = -4162
$ws.Range("A9:B14").Delete($xlShiftUp)
